$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 577.7778
$ws.Cells.Item(29, 12).Value = 2375.0001
$ws.Cells.Item(29, 9).Value = 150
$ws.Cells.Item(29, 11).Value = 450
$ws.Cells.Item(29, 10).Value = 791.6667
$ws.Cells.Item(29, 14).Value = -2937.0001
$ws.Cells.Item(29, 13).Value = -169
$ws.Cells.Item(31, 14).Value = -9460
$ws.Cells.Item(31, 8).Value = 2000
$ws.Cells.Item(31, 11).Value = 3000
$ws.Cells.Item(31, 9).Value = 1000
$ws.Cells.Item(31, 10).Value = 3000
$ws.Cells.Item(31, 12).Value = 9000
$ws.Cells.Item(31, 13).Value = -2770
$ws.Cells.Item(38, 9).Value = 132.11111
$ws.Cells.Item(38, 11).Value = 396.33333
$ws.Cells.Item(38, 8).Value = 2934.45
$ws.Cells.Item(38, 12).Value = 15681.819
$ws.Cells.Item(38, 14).Value = -16425.819
$ws.Cells.Item(38, 13).Value = -24.33332999999999
$ws.Cells.Item(38, 10).Value = 5227.273
$ws.Cells.Item(41, 11).Value = 323.66666
$ws.Cells.Item(41, 8).Value = 300.83334
$ws.Cells.Item(41, 13).Value = 116.33334
$ws.Cells.Item(41, 9).Value = 323.66666
$ws.Cells.Item(41, 12).Value = 278
$ws.Cells.Item(41, 10).Value = 278
$ws.Cells.Item(41, 14).Value = -1158
$ws.Cells.Item(42, 8).Value = 38
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = 116
$ws.Cells.Item(42, 11).Value = 114
$ws.Cells.Item(42, 9).Value = 38
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 14).ClearContents()
$ws.Cells.Item(43, 13).Value = -831
$ws.Cells.Item(43, 8).Value = 1084.7391
$ws.Cells.Item(43, 12).Value = 1183.2667
$ws.Cells.Item(43, 14).Value = -1321.2667
$ws.Cells.Item(43, 9).Value = 900
$ws.Cells.Item(43, 10).Value = 1183.2667
$ws.Cells.Item(43, 11).Value = 900

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(17, 11).Value = 5000
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 9).Value = 5000
$ws.Cells.Item(17, 8).Value = 5000
$ws.Cells.Item(17, 14).ClearContents()
$ws.Cells.Item(17, 13).Value = -4827
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 14).ClearContents()
$ws.Cells.Item(32, 14).Value = -107381.8
$ws.Cells.Item(32, 8).Value = 17643.744
$ws.Cells.Item(32, 9).Value = 3711.8594
$ws.Cells.Item(32, 11).Value = 3711.8594
$ws.Cells.Item(32, 12).Value = 106807.8
$ws.Cells.Item(32, 10).Value = 106807.8
$ws.Cells.Item(32, 13).Value = -3424.8594
$ws.Cells.Item(49, 8).Value = 12000
$ws.Cells.Item(49, 10).Value = 12000
$ws.Cells.Item(49, 12).Value = 12000
$ws.Cells.Item(49, 14).Value = -12520

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 10000
$ws.Cells.Item(18, 8).Value = 10000
$ws.Cells.Item(18, 12).Value = 10000
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 13).ClearContents()
$ws.Cells.Item(18, 14).Value = -11058
$ws.Cells.Item(20, 9).Value = 15877041
$ws.Cells.Item(20, 12).Value = 2256.182
$ws.Cells.Item(20, 14).Value = -2750.182
$ws.Cells.Item(20, 11).Value = 15877041
$ws.Cells.Item(20, 10).Value = 2256.182
$ws.Cells.Item(20, 13).Value = -15876794
$ws.Cells.Item(20, 8).Value = 10420084

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 430703.75
$ws.Cells.Item(5, 11).Value = 1328.47824
$ws.Cells.Item(5, 13).Value = -1216.47824
$ws.Cells.Item(5, 9).Value = 442.82608
$ws.Cells.Item(17, 11).Value = 399.99999
$ws.Cells.Item(17, 12).Value = 106299.996
$ws.Cells.Item(17, 9).Value = 133.33333
$ws.Cells.Item(17, 8).Value = 23666.666
$ws.Cells.Item(17, 14).Value = -106637.996
$ws.Cells.Item(17, 13).Value = -230.99999
$ws.Cells.Item(17, 10).Value = 35433.332
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 14).ClearContents()
$ws.Cells.Item(23, 11).Value = 141
$ws.Cells.Item(23, 14).Value = -629.999996
$ws.Cells.Item(23, 13).Value = 94
$ws.Cells.Item(23, 9).Value = 47
$ws.Cells.Item(23, 12).Value = 159.999996
$ws.Cells.Item(23, 10).Value = 53.333332
$ws.Cells.Item(23, 8).Value = 50.166668
$ws.Cells.Item(25, 9).Value = 625
$ws.Cells.Item(25, 13).Value = -1706
$ws.Cells.Item(25, 8).Value = 2450.1
$ws.Cells.Item(25, 11).Value = 1875
$ws.Cells.Item(30, 8).Value = 2450.1
$ws.Cells.Item(30, 11).Value = 1875
$ws.Cells.Item(30, 9).Value = 625
$ws.Cells.Item(30, 13).Value = -1773
$ws.Cells.Item(35, 12).Value = 6929.667
$ws.Cells.Item(35, 14).Value = -7505.667
$ws.Cells.Item(35, 11).Value = 3000
$ws.Cells.Item(35, 8).Value = 2178.9
$ws.Cells.Item(35, 10).Value = 2309.889
$ws.Cells.Item(35, 13).Value = -2712
$ws.Cells.Item(35, 9).Value = 1000
$ws.Cells.Item(36, 8).Value = 3058.3333
$ws.Cells.Item(36, 12).Value = 13200
$ws.Cells.Item(36, 14).Value = -13538
$ws.Cells.Item(36, 13).Value = -956
$ws.Cells.Item(36, 9).Value = 375
$ws.Cells.Item(36, 10).Value = 4400
$ws.Cells.Item(36, 11).Value = 1125
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 8).Value = 1840
$ws.Cells.Item(41, 13).ClearContents()
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 12).Value = 5520
$ws.Cells.Item(41, 10).Value = 1840
$ws.Cells.Item(41, 14).Value = -6196
$ws.Cells.Item(42, 8).Value = 2999.6667
$ws.Cells.Item(42, 12).Value = 8999.000100000001
$ws.Cells.Item(42, 10).Value = 2999.6667
$ws.Cells.Item(42, 14).Value = -10067.0001
$ws.Cells.Item(43, 8).Value = 6882.353
$ws.Cells.Item(43, 12).Value = 20647.059
$ws.Cells.Item(43, 14).Value = -20875.059
$ws.Cells.Item(43, 10).Value = 6882.353
$ws.Cells.Item(47, 8).Value = 668.625
$ws.Cells.Item(47, 9).Value = 621.2857
$ws.Cells.Item(47, 10).Value = 1000
$ws.Cells.Item(47, 11).Value = 1863.8571
$ws.Cells.Item(47, 12).Value = 3000
$ws.Cells.Item(47, 13).Value = -1432.8571
$ws.Cells.Item(47, 14).Value = -3862
$ws.Cells.Item(48, 10).Value = 2351
$ws.Cells.Item(48, 8).Value = 2145.3333
$ws.Cells.Item(48, 14).Value = -7553
$ws.Cells.Item(48, 12).Value = 7053
$ws.Cells.Item(48, 9).Value = 500
$ws.Cells.Item(48, 13).Value = -1250
$ws.Cells.Item(48, 11).Value = 1500
$ws.Cells.Item(49, 9).Value = 0
$ws.Cells.Item(49, 13).ClearContents()
$ws.Cells.Item(49, 8).Value = 7187.25
$ws.Cells.Item(49, 10).Value = 7187.25
$ws.Cells.Item(49, 12).Value = 21561.75
$ws.Cells.Item(49, 14).Value = -21873.75
$ws.Cells.Item(49, 11).Value = 0
$ws.Cells.Item(135, 9).Value = 442.82608
$ws.Cells.Item(135, 11).Value = 3985.43472
$ws.Cells.Item(135, 8).Value = 430703.75
$ws.Cells.Item(135, 13).Value = -1450.43472

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 14).Value = -3054.5
$ws.Cells.Item(24, 12).Value = 2708.5
$ws.Cells.Item(24, 8).Value = 7501354
$ws.Cells.Item(24, 10).Value = 2708.5
$ws.Cells.Item(47, 8).Value = 27666.666
$ws.Cells.Item(47, 10).Value = 27666.666
$ws.Cells.Item(47, 12).Value = 27666.666
$ws.Cells.Item(47, 14).Value = -28802.666
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 14).ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(42, 8).Value = 9933.333000000001
$ws.Cells.Item(42, 12).Value = 9933.333000000001
$ws.Cells.Item(42, 10).Value = 9933.333000000001
$ws.Cells.Item(42, 14).Value = -11059.333
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 14).ClearContents()
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 14).ClearContents()
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(49, 8).Value = 9933.333000000001
$ws.Cells.Item(49, 10).Value = 9933.333000000001
$ws.Cells.Item(49, 12).Value = 9933.333000000001
$ws.Cells.Item(49, 14).Value = -10227.333
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 14).ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(47, 8).Value = 5850
$ws.Cells.Item(47, 10).Value = 5850
$ws.Cells.Item(47, 12).Value = 5850
$ws.Cells.Item(47, 14).Value = -6994
